# Regenerate test data to include credit card account refund/credit transactions.
# Updates the pivot-style summary values on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Dining Out
$ws.Range("B4").Value = -81801.64
$ws.Range("E4").Value = -38293.06
$ws.Range("H4").Value = -60928.02
$ws.Range("K4").Value = -181022.72

# Row 5 - Freelance Income
$ws.Range("C5").Value = 126416.96
$ws.Range("F5").Value = 46505.97
$ws.Range("I5").Value = 165690.46
$ws.Range("K5").Value = 338613.39

# Row 6 - Groceries
$ws.Range("B6").Value = -66810.42999999999
$ws.Range("E6").Value = -38819.41
$ws.Range("H6").Value = -73376.16
$ws.Range("K6").Value = -179006

# Row 7 - Interest Income
$ws.Range("C7").Value = 892750.23
$ws.Range("F7").Value = 521104.54
$ws.Range("I7").Value = 840546.9
$ws.Range("K7").Value = 2254401.67

# Row 8 - Pets
$ws.Range("B8").Value = -32211.18
$ws.Range("E8").Value = -15438.74
$ws.Range("H8").Value = -33069.4
$ws.Range("K8").Value = -80719.32000000001

# Row 9 - Pharmacy
$ws.Range("B9").Value = -82469.08
$ws.Range("E9").Value = -31167.05
$ws.Range("H9").Value = -76872.91
$ws.Range("K9").Value = -190509.04

# Row 10 - Rent
$ws.Range("B10").Value = -32168.74
$ws.Range("E10").Value = -16132.08
$ws.Range("H10").Value = -31931.7
$ws.Range("K10").Value = -80232.52

# Row 11 - Shopping
$ws.Range("B11").Value = -102411.74
$ws.Range("E11").Value = -54028.1
$ws.Range("H11").Value = -123307.01
$ws.Range("K11").Value = -279746.85

# Row 12 - Taxes
$ws.Range("B12").Value = -58558.86
$ws.Range("E12").Value = -27644.43
$ws.Range("H12").Value = -54895.93
$ws.Range("K12").Value = -141099.22

# Row 13 - Transfer From
$ws.Range("D13").Value = 384725.52
$ws.Range("G13").Value = 185434.74
$ws.Range("J13").Value = 395603.59
$ws.Range("K13").Value = 965763.85

# Row 14 - Transfer To
$ws.Range("D14").Value = -384725.52
$ws.Range("G14").Value = -185434.74
$ws.Range("J14").Value = -395603.59
$ws.Range("K14").Value = -965763.85

# Row 15 - Utilities
$ws.Range("B15").Value = -34618.25
$ws.Range("E15").Value = -15629.56
$ws.Range("H15").Value = -32754.63
$ws.Range("K15").Value = -83002.44

# Row 16 - Wages & Salary
$ws.Range("C16").Value = 157971.31
$ws.Range("F16").Value = 76187.89999999999
$ws.Range("I16").Value = 179165.33
$ws.Range("K16").Value = 413324.54

# Row 17 - Total
$ws.Range("B17").Value = -491049.92
$ws.Range("C17").Value = 1177138.5
$ws.Range("E17").Value = -237152.43
$ws.Range("F17").Value = 643798.41
$ws.Range("H17").Value = -487135.76
$ws.Range("I17").Value = 1185402.69
$ws.Range("K17").Value = 1791001.49
